# Applies the "Second Problem 2-Socks in the Dark" edit:
#   1. "Problem:" -> "Problem" / "1" / ":" (three separate runs)
#   2. Appends a whole new "Problem2" (socks-in-the-dark, 2 pairs) worked
#      problem after the first problem's solution, moving the _GoBack
#      bookmark down onto the new solution's answer paragraph.

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Find-ParagraphByText($doc, [string]$text) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.TrimEnd([char]13) -eq $text) {
            return $p
        }
    }
    return $null
}

# --- Change 1: split "Problem:" into three runs: "Problem", "1", ":" ---
$problemPara = Find-ParagraphByText $d "Problem:"
if ($problemPara -eq $null) {
    throw "Could not find 'Problem:' paragraph"
}
$problemPara.Range.InsertXML(@"
<w:p xmlns:w='$wNs'><w:r><w:t>Problem</w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>
"@)

# --- Change 2: add the "Problem2" socks section, with its own Breaking it
#     apart / Potential Solutions / Evaluate potential solution / Solution
#     write-up, and relocate the _GoBack bookmark onto its answer line ---
$startPara = Find-ParagraphByText $d "To pick the least and be guaranteed would be 4. "
if ($startPara -eq $null) {
    throw "Could not find the 'would be 4' answer paragraph"
}

$newBlock = @"
<w:p xmlns:w='$wNs'>
  <w:r><w:t xml:space='preserve'>To pick the least </w:t></w:r>
  <w:r><w:t xml:space='preserve'>and be guaranteed would be 4. </w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'>
  <w:r><w:t>Problem2:</w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:r><w:t xml:space='preserve'>There are 20 socks in a drawer; 10 black, 6 brown, and 4 white how many will I need to pick to 2 matching pairs in the dark; with the least number of selecting through the drawer? </w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'>
  <w:r><w:t>Breaking it apart:</w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:r><w:t xml:space='preserve'>You have a better chance of picking black socks compared to others; less chance of picking white. </w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'>
  <w:r><w:t>Potential Solutions:</w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:r><w:t xml:space='preserve'>Can pick </w:t></w:r>
  <w:r><w:t>6, 9, or 12</w:t></w:r>
  <w:r><w:t xml:space='preserve'> </w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'>
  <w:r><w:t>Evaluate potential solution:</w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:r><w:t>Picking 12</w:t></w:r>
  <w:r><w:t xml:space='preserve'> and have a bett</w:t></w:r>
  <w:r><w:t>er chance of getting 2</w:t></w:r>
  <w:r><w:t xml:space='preserve'> pair</w:t></w:r>
  <w:r><w:t>s. Can pick 9</w:t></w:r>
  <w:r><w:t>, get one of each color, and the extra would be a match.</w:t></w:r>
  <w:r><w:t xml:space='preserve'> Pick 6</w:t></w:r>
  <w:r><w:t xml:space='preserve'> have 50% chance of picking black. </w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>Solution:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:r><w:t xml:space='preserve'>To pick the least and be guaranteed </w:t></w:r>
  <w:r><w:t>would be 9</w:t></w:r>
  <w:bookmarkStart w:id='0' w:name='_GoBack'/>
  <w:bookmarkEnd w:id='0'/>
  <w:r><w:t xml:space='preserve'>. </w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'>
  <w:r><w:t xml:space='preserve'> </w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'/>
<w:p xmlns:w='$wNs'/>
"@

# Replace from the start of the "would be 4" answer paragraph through the
# end of the document (all the trailing empty paragraphs) with the updated
# paragraph plus the brand-new Problem2 section and relocated bookmark.
$tailRange = $d.Range($startPara.Range.Start, $d.Content.End)
$tailRange.InsertXML($newBlock)
